# Implement the 'Banishment' card (and its 'Cannibal' neighbour) in the Cards sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Card #40 "The Contract" is now finished -> flip Done? from N to Y ---
$ws.Range("K42").Value = "Y"

# --- Row 43: new card #41 "Banishment" ---
$ws.Range("E43").Value = "Banishment"
$ws.Range("F43").Value = "UTILITY"
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = "Shuffle one unit on the field into its player's deck."
$ws.Range("K43").Value = "Y"
# E43 would otherwise pick up column E's own default style; align it with its
# freshly-written neighbours instead of that unrelated column formatting.
$ws.Range("E43").Style = $ws.Range("F43").Style

# --- Row 44: new card #42 "Cannibal" ---
$ws.Range("E44").Value = "Cannibal"
$ws.Range("F44").Value = "MINION"
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 4
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = "Once per turn, you can kill one of your units, fully restore this unit's health."
$ws.Range("K44").Value = "N"
$ws.Range("E44").Style = $ws.Range("F44").Style

# --- Keep the sheet view / selection in sync with the last edited card row ---
$excel.ActiveWindow.TopLeftCell = $ws.Range("E16")
$ws.Range("K43").Select() | Out-Null
